# Update field definitions: add "e.g." examples + type suffix, add description
# of follow-up schedule, mark several field names with a trailing "*", and
# fix a field/definition mismatch for the DMARDs-related rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Field", "Definition"),
    @("institution", "university or organization of the uploading member (e.g., University of Calgary) : string"),
    @("age*", "age of the participant at the time of the scan (e.g., 45) : int"),
    @("sex_assigned_at_birth*", "sex of participant assigned at birth (e.g., M, F, Other) : string"),
    @("weight_kg*", "weight of the participant at the time of the scan in kilograms (e.g., 75.2) : float"),
    @("height_cm", "height of the participant at the time of the scan in centimeters (e.g., 185) : int"),
    @("disease_status", "current medical diagnosis or condition of the participant (e.g., OA) : string"),
    @("smoking_status", "participant's smoking habits (e.g., Y/N) : string"),
    @("current_treatment", "current medical treatments or therapies the participant is undergoing : string"),
    @("bone_alterting_meds", "indicates if the participant is taking medications known to alter bone metabolism (e.g., Y/N) : string"),
    @("conventional_DMARDS", "indicates if the participant is taking conventional Disease-Modifying Antirheumatic Drugs (e.g., Y/N) : string"),
    @("biological_DMARDS", "indicates if the participant is taking biological Disease-Modifying Antirheumatic Drugs (e.g., Y/N) : string"),
    @("steroid_use", "indicates if the participant is currently using or has recently used steroids (e.g., Y/N) : string"),
    @("motion_score", "a quantitative measure of motion artifact during the scan (e.g., a score from 0-10) : int"),
    @("description_of_hand_scanned", "a more detailed description if a hand scan was performed (e.g., dominant hand, specific fingers) : string"),
    @("study_id*", "unique id of the study that the participant is taking part in (e.g., 3300_SPECTRA) : string"),
    @("time_interval_between_scans", "description of follow up schedule for patients in the study (e.g., BL, FU3mo, FU6mo) : string"),
    @("groups", "participant's assigned group within the study (e.g., control, treatment, disease subtype) : string"),
    @("scan_date*", "date scan was recorded (YYYY-MM-DD) : date"),
    @("file_type", "file extension of image uploaded (e.g., .isq, .dcm) : string"),
    @("joint_scanned", "specific anatomical joint that was scanned (e.g., right knee, left wrist, lumbar spine) : string"),
    @("length_of_scan_region", "total size of the scan in milimiters in 3 dimensions [x_mm, y_mm, z_mm] : 1x3 float"),
    @("voxel_spacing", "voxel size in 3 dimensions in milimeters [x_mm, y_mm, z_mm] : 1x3 float"),
    @("composite_id", "id of scan comprised of participant age, study_id, and scan_date (this field is automatically generated) : string")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

$ws.Range("A20").Select()
